# 02132019.xlsx — "add jan and feb layers"
#
# Row 34 had a data-entry error: the location-name string that belongs in
# column E ("CHANNING/TELEGRAPH") had leaked into column F (shifting the
# latitude/longitude pair one column to the right, into G/H, and leaving no
# room for the trailing "NA" marker that every other row carries in column
# H). This fixes it so row 34 matches the layout of every other row:
#   F34 = latitude (number)
#   G34 = longitude (number)
#   H34 = "NA"      (string, same look as the other H-column "NA" cells)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the formatting of a cell that already has the correct "NA" style
# (style used throughout column H) and stamp it onto H34 before putting the
# text in, so H34 ends up s="2" t="s" just like its neighbors.
$ws.Range("E34").Copy()
$ws.Range("H34").PasteSpecial(-4122)  # xlPasteFormats

# Shift the latitude/longitude values left into F/G, and place the "NA"
# marker in H.
$ws.Range("F34").Value = 37.866916000000003
$ws.Range("G34").Value = -122.258786
$ws.Range("H34").Value = "NA"

# F34 previously held a text value styled like the other text columns;
# now that it holds a plain number again it should go back to the
# worksheet's default (unstyled) look, like the F column everywhere else.
$ws.Range("F34").ClearFormats()

# Selection marker left on H35 (matches the cursor position recorded the
# last time the sheet was saved).
[void]$ws.Range("H35").Select()
